$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 152
$ws.Range("A3").Value = 155
$ws.Range("A4").Value = 159
$ws.Range("A5").Value = 157
$ws.Range("A6").Value = 148
$ws.Range("A7").Value = 158
$ws.Range("A8").Value = 163
$ws.Range("A9").Value = 169
